# Insert a new data row at row 355, pushing the existing rows 355-454
# down to 356-455, and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 355 (shifts 355..454 -> 356..455)
$ws.Rows.Item(355).Insert()

# Populate the newly inserted row 355 with the new record's data
$ws.Cells.Item(355, 1).Value = 9
$ws.Cells.Item(355, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(355, 3).Value = "Metropolitana"
$ws.Cells.Item(355, 4).Value = 45211
$ws.Cells.Item(355, 5).Value = 13
$ws.Cells.Item(355, 6).Value = 100112021
$ws.Cells.Item(355, 7).Value = "Ají"
$ws.Cells.Item(355, 8).Value = "Inferno"
$ws.Cells.Item(355, 9).Value = "Primera"
$ws.Cells.Item(355, 10).Value = 124
$ws.Cells.Item(355, 11).Value = 29000
$ws.Cells.Item(355, 12).Value = 31000
$ws.Cells.Item(355, 13).Value = 30000
$ws.Cells.Item(355, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(355, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(355, 16).Value = 3000
$ws.Cells.Item(355, 17).Value = 10
$ws.Cells.Item(355, 18).Value = "Hortaliza"

# Match the date-formatted style used by column D in the surrounding rows
$ws.Cells.Item(355, 4).NumberFormat = $ws.Cells.Item(356, 4).NumberFormat
